# Swap the data blocks of rows 2-9 with rows 10-17 (columns B through F),
# leaving column A (message_id) untouched for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current values for columns B:F for rows 2-9 and rows 10-17
$topBlock = $ws.Range("B2:F9").Value2
$bottomBlock = $ws.Range("B10:F17").Value2

# Write them back swapped
$ws.Range("B2:F9").Value2 = $bottomBlock
$ws.Range("B10:F17").Value2 = $topBlock
